$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = "Jean Carlos"

# C2 looks like a date ("01/01/2021"); force it to stay plain text by using a
# leading apostrophe (quote-prefix), then reset the style back to Normal so
# no stray formatting (quotePrefix) sticks around on the cell.
$ws.Range("C2").Value = "'01/01/2021"
$ws.Range("C2").Style = "Normal"

# Remove the "Mensagem adicional" note for row 2 entirely.
$ws.Range("J2").ClearContents()

# --- Row 3 ---
$ws.Range("B3").Value = "Carlos "

$ws.Range("C3").Value = "'02/02/2021"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "10:45"

$ws.Range("E3").Value = 10

$ws.Range("J3").ClearContents()

# --- Row 4 ---
$ws.Range("B4").Value = "Naejc"

$ws.Range("C4").Value = "'01/01/2021"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").Value = "10:45"

$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 10

$ws.Range("J4").ClearContents()
